$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-05-30 Friday" "2025-05-31 Saturday"

Replace-Text "643÷5=128, 3" "237÷8=29, 5"
Replace-Text "914÷8=114, 2" "612÷5=122, 2"
Replace-Text "646÷5=129, 1" "972÷7=138, 6"
Replace-Text "306÷8=38, 2" "916÷4=229, 0"
Replace-Text "502÷8=62, 6" "555÷5=111, 0"

Replace-Text "882÷8=110, 2" "585÷4=146, 1"
Replace-Text "405÷9=45, 0" "984÷7=140, 4"
Replace-Text "416÷5=83, 1" "120÷2=60, 0"
Replace-Text "464÷7=66, 2" "455÷7=65, 0"
Replace-Text "906÷9=100, 6" "423÷2=211, 1"

Replace-Text "145÷9=16, 1" "939÷4=234, 3"
Replace-Text "751÷3=250, 1" "585÷6=97, 3"
Replace-Text "104÷2=52, 0" "734÷3=244, 2"
Replace-Text "981÷7=140, 1" "202÷7=28, 6"
Replace-Text "343÷7=49, 0" "490÷3=163, 1"

Replace-Text "591÷7=84, 3" "230÷5=46, 0"
Replace-Text "655÷7=93, 4" "313÷8=39, 1"
Replace-Text "417÷8=52, 1" "733÷4=183, 1"
Replace-Text "539÷4=134, 3" "929÷4=232, 1"
Replace-Text "460÷5=92, 0" "177÷5=35, 2"

Replace-Text "641÷3=213, 2" "259÷8=32, 3"
Replace-Text "190÷2=95, 0" "512÷6=85, 2"
Replace-Text "495÷5=99, 0" "772÷6=128, 4"
Replace-Text "258÷5=51, 3" "355÷2=177, 1"
Replace-Text "782÷5=156, 2" "356÷9=39, 5"

Write-Output "Done applying replacements"
